# Sprint 8 | maj
# Record today's progress on the burndown tracker:
#  - US#1 (row 2, "US#1 Vincenzo Eliott | BDD") logs 1h on each of the first two days
#  - US#3 (row 4, "US#3 Nico Max | Navigation v2") logs 1h on the second day

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("D4").Value = 1

# Move the active selection to A6, matching where the author left off editing.
$ws.Range("A6").Select() | Out-Null
